# Regenerate s_val data to filter save games: update B2:E11 and G2:G11
# with new computed values (F column / wins are unchanged).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @{ B = 3.272327238179451;  C = 1.626987699542094; D = 3.223369029078222;  E = 0.5333859586016987; G = 8.656069925401464 }
    3  = @{ B = 0.2881169905109251; C = 1.626987699542094; D = 0.7210945179870265; E = 0.5333859586016987; G = 3.169585166641744 }
    4  = @{ B = 1.445647641019636;  C = 1.626987699542094; D = 0.1496068669990043; E = 0.5333859586016987; G = 3.755628166162433 }
    5  = @{ B = 3.272327238179451;  C = 1.626987699542094; D = 0.1496068669990043; E = 0.5333859586016987; G = 5.582307763322248 }
    6  = @{ B = 3.272327238179451;  C = 1.626987699542094; D = 0.7210945179870265; E = 0.5333859586016987; G = 6.15379541431027 }
    7  = @{ B = 1.445647641019636;  C = 1.626987699542094; D = 3.223369029078222;  E = 0.5333859586016987; G = 6.82939032824165 }
    8  = @{ B = 3.272327238179451;  C = 1.626987699542094; D = 0.7210945179870265; E = 0.5333859586016987; G = 6.15379541431027 }
    9  = @{ B = 3.272327238179451;  C = 1.626987699542094; D = 3.223369029078222;  E = 0.5333859586016987; G = 8.656069925401464 }
    10 = @{ B = 1.445647641019636;  C = 1.626987699542094; D = 3.223369029078222;  E = 0.5333859586016987; G = 6.82939032824165 }
    11 = @{ B = 1.445647641019636;  C = 1.626987699542094; D = 0.1496068669990043; E = 0.5333859586016987; G = 3.755628166162433 }
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Range("B$row").Value = $vals.B
    $ws.Range("C$row").Value = $vals.C
    $ws.Range("D$row").Value = $vals.D
    $ws.Range("E$row").Value = $vals.E
    $ws.Range("G$row").Value = $vals.G
}
